$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text (inline string) type rather than
# being auto-converted to numbers by Excel when the new value looks numeric.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '68.755.32'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +2.65%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.658.03'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +2.48%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '600.82'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +2.04%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '155.60'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +4.32%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.66%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.656.01'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +2.46%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.138'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +12.87%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.23'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +2.00%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.05'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +3.56%  '
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +6.12%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.143.56'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '68.585.17'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +2.30%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.657.59'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +2.54%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.45'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +4.49%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '367.12'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +1.35%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.44'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +2.01%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.28'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.88'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.31%  '
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +5.13%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '72.65'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.02'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +1.25%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +8.56%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.787.98'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +2.82%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.21%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '576.82'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.46%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +4.89%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.98'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +5.40%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +3.55%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +5.40%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +4.09%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '159.38'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +2.00%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +5.48%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '19.29'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +2.26%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.42'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +4.96%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.368'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +1.18%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.66'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +7.85%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0₆0322'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +14.47%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '156.54'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +3.21%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.74'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.11%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.72'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +2.94%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '22.03'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +3.98%  '
